$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 16 play-by-play yardage logs ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = '1 0 5 5 9 4 5 4 6 3 13 6 2 7 13 3 -1 2 9 0 1 3 2 15 3 3 3 9 1 4 3 -1 3 7 1 5 -4 5 4 8 5 1 -1 4 5 14 1 2 2 2 6 0 25 9 4 0 13 3 5 -2 6 0 -2 -1 4 2 2 6 -1 3 21 1 2 3 2 2 18 2 3 9 9 -2 1 2 14 2 2 1 3 18 5 2 -4 11 0 22 1 1 2 9 7 4 4 -3 7 0 4 3 2 13 6 4 8 6 2 -1 17 3 8 2 15 9 2 1 0 2 -1 3 10 11 3 4 3 1 3 1 12 10 10 -1 8 0 1 6 1 4 1 8 0 1 3 1 11 7 -1 5 -4 2 1 2 0 7 1 3 5 4 2 5 7 17 12 -4 5 6 2 9 2 3 2 8 14 5 3 9 12 4 4 6 5 4 7 5 3 6 13 4 6 0 10 3 -1 1 2 -4 4 4 1 6 1 5 6 4 7 6 1 7 6 -1 4 1 10 -3 8 1 3 5 3 2 6 1 2 1 4 2 1 9 4 3 0 5 4 7 3 0 4 21 9 5 3 4 3 0 6 2 3 2 10 6 5 13 10 3 13 6 2 3 -4 3 -2 4 3 -1 1 -2 8 3 4 1 9 5 -4 4 2 1 4 11 14 -2 2 6 10 3 9 5 3 7 12 5 6 4 7 10 2 10 6 -1 2 12 3 5 3 3 5 2 2 0 3 9 2 1 3 13 4 6 0 5 1 10 1 7 6 4 0 4 1 7 2 3 7 5 45 1 3 0 3 3 1 7 3 7 0 1 -2 -2 -1 3 1 -3 3 -2 0 4 7 1 0 3 0 5 6 4 6 5 8 6 10 3 7 8 4 1 7 2 3 1 5 8 4 3 2 2 9 0 14 3 -6 7 15 0 18 3 0 3 5 3 2 6 1 1 0 2 5 -1 -1 16 6 0 7 3 5 9 0 2 -2 0 11 0 -1 1 3 2 -5 3 9 0 3 3 5 -3 10 4 9 1 2 1 0 4 13 7 2 2 -1 4 4 12 0 2 14 0 1 6 6 4 14 12 2 -1 1 -1 14 7 4 11 -2 2 6 7 3 3 0 7 -2 7 7 5 1 5 6 0 26 3 1 5 0 3 3 -4 1 2 9 -5 8 2 3 11 1 8 6 11 3 5 -1 2 5 3 2 8 5 3 6 2 3 5 -2 6 9 6 0 5 2 5 2 6 30 0 13 3 6 13 7 -3 2 -2 2 2 4 7 -4 0 1 3 -2 2 4 9 3 9 2 2 1 -2 6 2 5 2 7 2 13 6 6 3 3 6 20 9 2 2 6 4 6 6 6 3 6 6 -2 1 11 8 2 6 -1 1 5 3 2 1 2 4 2 0 4 4 4 -2 3 10 2 1 3 9 6 0 3 0 4 10 1 4 4 1 -1 2 11 11 12 10 6 2 4 12 1 16 15 10 2 5 4 1 1 4 0 5 2 10 4 -9 3 15 4 8 6 -4 -3 1 7 -2 8 6 3 0 -3 4 3 8 3 4 9 6 24 2 4 0 3 6 10 4 4 4 2 1 15 5 2 2 3 1 1 -3 0 5 -1 5 3 -1 11 8 -1 12 8 4 2 4 1 4 1 4 -5 -1 5 11 17 4 0 1 2 -2 2 4 1 11 4 21 3 -2 0 7 8 -4 5 2 2 0 2 30 3 4 1 13 -1 1 7 1 7 8 5 7 5 33 2 -2 1 6 2 2 -2 0 4 5 0 3 3 4 11'
$ws.Range("C2").Value = '5 1 7 1 5 1 4 4 7 14 5 -1 7 12 15 2 -4 7 2 7 9 4 -2 0 2 3 6 1 2 5 0 4 1 7 2 5 4 4 1 0 4 4 1 -5 6 -1 25 6 4 3 46 0 8 10 4 4 9 6 -1 12 3 0 2 7 4 3 1 0 8 -1 4 4 11 6 9 1 4 3 7 4 3 0 3 3 -2 5 -3 8 48 2 8 4 7 9 7 -1 11 8 7 35 4 4 4 7 24 14 1 3 0 5 4 2 7 11 6 3 1 2 4 2 3 1 3 4 -3 3 -1 7 1 0 -1 5 11 0 0 1 12 0 9 2 7 1 4 0 1 18 1 7 1 6 25 4 1 1 32 0 0 4 7 2 4 2 0 2 4 -4 6 6 3 2 1 1 4 2 0 12 8 4 3 1 1 13 9 0 2 4 -4 11 2 2 1 3 10 14 4 1 6 8 3 2 1 3 -4 5 2 0 2 7 4 8 0 3 0 1 8 3 0 7 0 3 1 14 4 2 -3 1 6 3 11 1 98 1 -3 3 -2 1 -1 4 5 4 12 1 6 6 2 3 8 4 3 2 3 0 4 1 6 4 0 2 8 -1 1 4 0 0 4 6 2 9 0 6 2 6 -1 2 5 11 5 2 4 4 10 1 14 3 0 9 3 3 2 5 0 4 6 5 9 1 4 4 9 -1 0 15 0 2 11 5 4 -1 6 1 -4 46 6 2 6 3 8 14 4 6 5 7 3 8 5 8 18 10 3 2 10 2 5 6 3 7 -1 2 5 7 5 10 6 13 3 4 0 22 9 3 -1 -6 3 4 7 12 1 3 4 -2 4 6 1 36 3 4 10 5 3 3 6 15 3 14 1 12 7 0 0 -2 6 5 4 2 2 1 1 2 6 2 2 0 2 2 1 1 1 9 2 12 2 2 0 0 4 2 12 0 1 7 1 5 9 -7 8 6 5 6 6 4 2 -1 1 -2 7 1 4 3 1 1 3 1 0 2 2 19 0 8 7 5 7 1 4 0 1 2 6 1 11 5 5 47 14 0 14 0 0 10 5 2 18 8 1 5 6 4 1 3 0 1 4 4 1 4 14 3 6 1 6 7 18 -3 12 1 -3 2 -3 3 4 14 5 9 1 3 4 3 5 1 8 9 6 0 14 2 0 30 16 -1 -1 23 4 0 3 3 -1 -3 16 5 4 16 14 -1 7 2 2 0 1 9 3 0 -2 -2 7 11 3 4 2 3 6 2 2 -2 4 19 0 4 9 2 4 1 -1 9 0 5 2 18 2 8 -2 3 3 14 5 0 5 3 1 2 2 0 16 12 1 5 2 -1 6 2 13 3 3 2 8 4 11 12 2 -2 3 6 5 1 4 11 2 5 1 7 -6 1 2 9 4 0 -1 1 6 2 4 0 4 0 0 4 4 6 11 11 1 9 4 3 10 2 -2 4 -3 15 5 2 11 1 8 3 4 0 1 4 5 11 20 2 0 11 9 4 9 4 2 8 -2 0 14 3 11 0 0 7 0 3 2 6 3 1 5 0 2 2 3 2 5 -3 3 0 2 16 0 2 4 4 3 3 8 7 2 0 8 2 1 9 2 1 1 5 -1 0 2 6 5 4 11 0 4 -1 4 11 4 0 11 1 4 3 11 8 5 9 0 2 1 -2 1 11 11 4 2 3 5 3 8 2 16 6 3 0 2 2 26 3 4 1 -2 12 3 1 2 4 1 8 7 4 5 2 1 5 3 4 5 55 0 4 4 4 -2 9 1 2 3 -1 3 1 1 2 0 11 14 7 2 -3 16 -1 12 -1 0 -1'
$ws.Range("B3").Value = '10 2 16 5 9 16 4 11 7 5 11 11 11 0 5 18 6 1 11 16 75 19 6 10 19 -1 8 9 7 13 12 3 3 7 4 15 5 10 11 30 39 16 18 8 7 5 12 12 2 22 14 8 13 7 14 12 18 -1 9 24 2 9 14 13 8 13 4 10 4 4 38 1 6 5 6 27 9 14 3 3 17 12 23 24 6 11 10 -1 16 17 3 20 24 20 -3 2 16 8 3 5 8 15 0 12 10 8 13 7 7 5 17 -1 29 14 1 12 57 2 8 12 23 22 3 1 17 14 3 -1 9 9 22 12 5 5 16 1 20 9 38 2 13 19 -4 10 2 22 74 5 11 19 8 4 8 7 8 3 4 26 2 5 11 7 6 8 5 -1 5 11 27 3 29 21 4 2 4 6 -7 5 -1 13 5 42 7 4 1 12 7 0 5 11 7 14 9 3 11 12 -6 8 8 3 28 24 3 11 2 9 22 11 5 3 9 6 23 0 5 0 2 5 -3 7 6 38 24 7 5 1 5 17 2 2 4 17 9 12 8 14 52 10 2 4 6 3 8 7 15 7 18 20 17 -3 11 12 10 -2 7 11 4 17 18 1 5 41 7 25 11 6 10 5 34 9 13 35 4 12 0 5 13 4 8 3 3 8 18 4 0 10 21 12 0 9 19 17 1 2 14 8 23 16 0 4 20 32 8 2 37 15 8 5 16 20 7 -6 6 18 7 9 6 9 15 13 5 14 40 11 9 8 1 16 -1 1 44 0 14 13 10 12 14 14 10 4 22 3 3 4 -1 8 14 25 9 8 12 33 9 3 8 44 23 31 22 18 11 13 -5 8 1 9 26 5 27 15 11 57 14 4 5 11 22 17 1 7 12 7 3 -2 32 20 20 19 8 5 13 18 19 5 16 7 17 2 10 6 7 13 8 1 4 8 11 11 16 9 7 1 29 12 26 1 13 8 24 12 7 8 10 10 16 25 10 8 30 5 8 11 23 9 29 9 3 13 6 5 17 1 8 5 9 39 9 15 2 11 4 7 13 8 6 5 5 8 6 13 14 5 6 12 11 -1 8 5 7 0 7 5 24 2 8 9 4 7 25 9 0 16 18 17 5 15 14 1 4 4 21 6 4 7 4 2 9 1 3 5 5 12 21 5 12 6 13 7 5 1 7 7 19 4 17 -2 13 10 23 7 9 21 4 11 1 11 13 6 2 6 5 33 13 -3 6 21 19 17 6 9 7 12 13 17 2 6 16 -2 6 10 2 3 11 6 9 2 10 6 10 1 2 3 9 8 15 7 3 8 13 10 5 3 11 2 19 7 9 4 10 16 27 7 6 4 8 11 7 64 4 13 4 7 18 24 13 17 7 11 15 22 10 23 9 8 14 6 18 7 9 3 19 19 13 18 11 1 9 8 8 9 4 11 13 9 16 4 23 5 6 8 5 8 14 8 10 5 15 63 0 4 6 2 9 19 8 6 13 2 18 14 6 6 17 6'
$ws.Range("C3").Value = '7 8 45 8 6 10 11 3 2 4 6 8 23 9 -1 15 29 0 18 6 14 8 6 7 2 50 6 23 2 7 4 6 12 8 36 12 11 2 0 5 8 -1 4 16 28 6 11 11 9 9 10 6 9 3 17 5 11 8 12 1 22 -1 16 -3 7 6 7 6 10 11 14 7 5 15 2 11 21 2 7 -3 8 3 8 7 3 2 6 8 -5 8 9 -2 11 25 -3 13 3 1 12 8 2 4 2 13 7 8 -1 20 6 9 6 8 8 17 4 26 13 26 11 9 8 1 9 17 11 2 7 4 2 8 12 11 7 10 13 39 6 13 6 8 8 18 6 7 1 23 7 6 9 9 2 7 4 15 8 14 10 5 17 3 18 7 22 17 4 6 16 6 14 23 -2 5 7 20 0 8 7 24 28 11 7 13 20 16 18 8 19 19 7 7 23 6 0 12 9 9 5 20 14 10 3 8 6 32 23 7 5 8 14 1 10 9 10 28 14 29 -5 3 44 4 10 26 28 8 13 2 3 6 23 19 14 3 5 3 9 8 7 2 4 15 16 44 3 6 6 19 22 6 8 3 25 13 7 7 6 21 9 11 5 1 4 19 5 5 35 8 14 18 3 3 5 12 6 5 2 12 3 8 -3 6 9 7 12 6 5 9 17 9 16 0 7 15 18 6 17 5 6 17 10 15 6 7 12 25 10 10 3 8 4 13 7 6 2 3 4 2 31 2 11 37 14 32 29 2 20 49 1 9 12 7 1 1 14 2 14 6 12 4 7 14 6 4 -2 10 22 2 2 5 -2 4 17 12 0 8 50 14 3 28 7 5 3 8 4 19 7 4 5 9 23 7 13 29 8 5 5 18 7 8 5 5 9 8 10 17 9 27 18 2 4 9 6 11 4 6 7 35 10 14 10 8 -3 1 11 20 7 8 22 6 15 2 40 12 8 9 7 4 16 3 7 4 27 3 8 23 -1 11 4 5 8 30 10 11 1 10 5 5 13 15 11 5 12 5 18 12 10 14 18 4 7 20 35 6 7 6 23 8 -7 6 -5 3 10 8 8 2 7 2 0 1 15 17 7 11 8 13 5 53 14 20 11 9 6 7 7 3 3 20 7 24 13 2 -6 8 14 7 9 24 35 5 17 26 40 3 4 12 6 8 6 7 13 -3 27 15 -2 6 18 7 1 27 15 7 16 15 1 6 8 15 -4 2 14 16 2 5 8 9 12 4 12 3 7 7 15 11 2 2 18 1 -3 2 9 17 5 8 15 10 7 8 12 -1 4 9 41 3 28 3 7 18 1 4 19 5 20 4 0 5 19 21 7 9 -5 17 5 3 30 7 39 6 4 5 12 12 7 13 4 15 10 18 2 5 -1 25 4 10 6 0 9 57 9 5 5 5 7 11 11 7 7 4 11 11 9 2 -2 5 6 -2 6 36 9 13 3 20 4 4 14 12 7 14 3 11 6 2 23 12 9 15 9 10 12 11 4 10 9 7 11 7 7 7 8 28 20 4 3 4 13 7 19 7 1 10 22 62 13 23 11 12 7 10'

# --- OFF sheet ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 400
$ws.Range("F2").Value = 90
$ws.Range("G2").Value = 117
$ws.Range("J2").Value = 65
$ws.Range("L2").Value = 555
$ws.Range("M2").Value = 336
$ws.Range("O2").Value = 57
$ws.Range("Q2").Value = 991
$ws.Range("C3").Value = 330
$ws.Range("E3").Value = 71
$ws.Range("F3").Value = 229
$ws.Range("H3").Value = 49
$ws.Range("I3").Value = 126
$ws.Range("J3").Value = 89
$ws.Range("N3").Value = 38

# --- DEF sheet ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 395
$ws.Range("D2").Value = 19
$ws.Range("F2").Value = 108
$ws.Range("G2").Value = 102
$ws.Range("J2").Value = 58
$ws.Range("L2").Value = 502
$ws.Range("M2").Value = 331
$ws.Range("Q2").Value = 987
$ws.Range("B3").Value = 14
$ws.Range("C3").Value = 302
$ws.Range("D3").Value = 13
$ws.Range("E3").Value = 69
$ws.Range("F3").Value = 218
$ws.Range("G3").Value = 55
$ws.Range("H3").Value = 61
$ws.Range("I3").Value = 100
$ws.Range("J3").Value = 107

# --- ST sheet ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 145
$ws.Range("D2").Value = 106
$ws.Range("J2").Value = 55
$ws.Range("K2").Value = 53
$ws.Range("B3").Value = 114
$ws.Range("D3").Value = '60 60 44 57 48 50 52 42 40 40 54 51 33 45 21 44 43 67 49 58 43 53 53 40 44 35 56 33 42 43 50 45 49 39 48 33 44 67 39 51 49 36 18 54 35 46 36 33 39 51 37 32 42 37 49 49 43 44 42 44 30 40 33 40 40 44 39 44 36 28 39 46 32 40 34 39 49 36 46 44 57 59 46 43 44 55 54 48 50 32 55 57 55 40 51 54 33 47 54 41 40 40 50 55 48 60'
$ws.Range("B6").Value = '38 21 24 27 35 24 18 23 17 24 21 19 35 16 20 26 23 23 28 98 25 9 16 23 28 21 29 21 21 21 39 13 31 24 16 13 10 21 22 26 21 26 28 0 25 21 10 18 21 27 25 26'
$ws.Range("D4").Value = '27 10 0 0 0 0 10 0 0 0 3 4 0 0 0 6 0 0 0 0 0 0 12 4 0 0 83 0 0 0 0 0 5 0 0 0 0 0 10 0 10 0 0 15 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 22 0 0 7 0 0 0 0 0 0 2 12 6 0 4 8 0 13 0 18 15 7 0 1 17 7 0 0 0 0 7 7 0 0 0 0 17 14 0 0 4 9 20 6 5'
$ws.Range("D5").Value = '19 10 0 0 -1 0 8 0 0 0 0 6 0 0 0 -3 0 0 0 9 0 0 0 8 0 0 0 7 14 0 11 20 0 0 0 13 0 0 0 1 19 11 0 0 0 0 8 0 0 0 4 9 0 0 0 12 15 0 10 14 0 0 0 5 0 12 0 0 0 0 0 0 0 0 0 22 8 0 9 1 0 4 0 0 10 0 0 4 0 0 0 0 0 0 4 0 10 0 16 10 0 0 0 0 0 0 0 0 0 0 0 0 16 16 6 0'

# --- TURNS sheet ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("B2").Value = 23
$ws.Range("D2").Value = 13
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 17

# --- PEN sheet ---
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("D2").Value = 21
$ws.Range("D3").Value = 12
